$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: "Objetivos:" B/C value replaced with faculty info
$ws.Range("B10").Value = '8855158 - Morun Bernardino Neto'
$ws.Range("C10").Value = '8855158 - Morun Bernardino Neto'

# Drop old rows 13-24 entirely, then rebuild rows 13-23 with the new layout/heights
$ws.Rows("13:24").Delete()

# Row 13
$ws.Range("A13").Value = 'Programa resumido:'
$ws.Range("B13").Value = 'Semestral'
$ws.Range("C13").Value = 'Semestral'
$ws.Rows(13).RowHeight = 60

# Row 14
$ws.Range("A14").Value = 'Short syllabus:'
$ws.Range("B14").Value = 'Definition of epidemiology and its areas of activity; Types of epidemiological studies; Errors and confounding factors in epidemiological studies; Environmental epidemiology: exposure and quantification of exposure to environmental factors; Risk assessment; Impact assessment on the environment and public health; Risk management.'
$ws.Range("C14").Value = 'Definition of epidemiology and its areas of activity; Types of epidemiological studies; Errors and confounding factors in epidemiological studies; Environmental epidemiology: exposure and quantification of exposure to environmental factors; Risk assessment; Impact assessment on the environment and public health; Risk management.'
$ws.Rows(14).RowHeight = 60

# Row 15
$ws.Range("A15").Value = 'Programa:'
$ws.Range("B15").Value = '01/01/2022'
$ws.Range("C15").Value = '01/01/2022'
$ws.Rows(15).RowHeight = 120

# Row 16
$ws.Range("A16").Value = 'Syllabus:'
$ws.Range("B16").Value = 'Introduction: Definition, field of action of epidemiology; Epidemiology and public health. Types of studies: Observational epidemiological studies (descriptive studies, ecological or correlation studies, ecological fallacy, cross-sectional studies, case and control studies, cohort studies); Experimental epidemiological studies (randomized clinical trial, field trials, community trials). Errors and confounding factors: Potential errors in epidemiological studies (random error, sample size, systematic error, selection bias, measurement bias); Confounding factors (control of confounding factors, validity, ethical issues). Environmental epidemiology: Exposure to environmental factors and quantification of exposure: Biological monitoring; Interpretation of biological data; Individual measures versus group measures; Population dose; Dose-effect relationship and dose-response relationship. Risk: Risk assessment; Health impact assessment; Risk management; Environmental and public health impact assessment.'
$ws.Range("C16").Value = 'Introduction: Definition, field of action of epidemiology; Epidemiology and public health. Types of studies: Observational epidemiological studies (descriptive studies, ecological or correlation studies, ecological fallacy, cross-sectional studies, case and control studies, cohort studies); Experimental epidemiological studies (randomized clinical trial, field trials, community trials). Errors and confounding factors: Potential errors in epidemiological studies (random error, sample size, systematic error, selection bias, measurement bias); Confounding factors (control of confounding factors, validity, ethical issues). Environmental epidemiology: Exposure to environmental factors and quantification of exposure: Biological monitoring; Interpretation of biological data; Individual measures versus group measures; Population dose; Dose-effect relationship and dose-response relationship. Risk: Risk assessment; Health impact assessment; Risk management; Environmental and public health impact assessment.'
$ws.Rows(16).RowHeight = 120

# Row 17
$ws.Range("A17").Value = 'Avaliação:'

# Row 18
$ws.Range("A18").Value = 'Método:'
$ws.Range("B18").Value = '8855158 - Morun Bernardino Neto'
$ws.Range("C18").Value = '8855158 - Morun Bernardino Neto'
$ws.Rows(18).RowHeight = 60

# Row 19
$ws.Range("A19").Value = 'Critério:'
$ws.Range("B19").Value = 'Aulas teóricas expositivas com resolução de exercícios e discussão de casos reais de impactos ambientais e seus potenciais reflexos à saúde pública: análise de riscos, avaliação dos impactos ambientais, avaliação dos impactos à saúde pública e manejo de riscos.'
$ws.Range("C19").Value = 'Aulas teóricas expositivas com resolução de exercícios e discussão de casos reais de impactos ambientais e seus potenciais reflexos à saúde pública: análise de riscos, avaliação dos impactos ambientais, avaliação dos impactos à saúde pública e manejo de riscos.'
$ws.Rows(19).RowHeight = 60

# Row 20
$ws.Range("A20").Value = 'Norma de recuperação:'
$ws.Range("B20").Value = 'O sistema de avaliação será composto por 2 avaliações de igual peso. A Nota Final será obtida por meio da média simples dessas duas avaliações. Estará aprovado por notas o aluno que obtiver nota final igual ou superior a 5,0 pontos.(Nota final+P_recuperação)/2'
$ws.Range("C20").Value = 'O sistema de avaliação será composto por 2 avaliações de igual peso. A Nota Final será obtida por meio da média simples dessas duas avaliações. Estará aprovado por notas o aluno que obtiver nota final igual ou superior a 5,0 pontos.(Nota final+P_recuperação)/2'
$ws.Rows(20).RowHeight = 60

# Row 21
$ws.Range("A21").Value = 'Bibliografia:'
$ws.Range("B21").Value = 'Estará em período de recuperação o aluno que obtiver notas entre 3,0 e 4,9. Para esses alunos a Nota Finalrec será calculada pela média simples entre a avaliação de recuperação (todo o conteúdo do semestre) e sua nota final.(Nota final+P_recuperação)/2'
$ws.Range("C21").Value = 'Estará em período de recuperação o aluno que obtiver notas entre 3,0 e 4,9. Para esses alunos a Nota Finalrec será calculada pela média simples entre a avaliação de recuperação (todo o conteúdo do semestre) e sua nota final.(Nota final+P_recuperação)/2'
$ws.Rows(21).RowHeight = 120

# Row 22
$ws.Range("A22").Value = 'Requisitos:'

# Row 23
$ws.Range("B23").Value = 'LOB1012 -  Estatística  (Requisito fraco)
'
$ws.Range("C23").Value = 'LOB1012 -  Estatística  (Requisito fraco)
'
$ws.Rows(23).RowHeight = 30

# The sheets <cols> definitions have an overlapping range for column B
# (min=1,max=2 style=1 AND min=2,max=2 style=2), so brand-new column-B cells
# default to the wrong (bold) style. Re-apply the correct body style (from an
# existing, correctly styled column-B cell) to every newly written B-cell.
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B14").PasteSpecial(-4122) | Out-Null
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B15").PasteSpecial(-4122) | Out-Null
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B16").PasteSpecial(-4122) | Out-Null
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B18").PasteSpecial(-4122) | Out-Null
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B19").PasteSpecial(-4122) | Out-Null
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B20").PasteSpecial(-4122) | Out-Null
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B21").PasteSpecial(-4122) | Out-Null
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B23").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false